$wb = $excel.ActiveWorkbook

# --- LogIn sheet: update the sample login data ---
$ws = $wb.Worksheets.Item("LogIn")

$ws.Range("A2").Value = 1005254554
$ws.Range("B2").Value = "7474"

# Make LogIn the active sheet and select A2, matching the new active tab / selection.
$ws.Activate()
[void]$ws.Range("A2").Select()
